# Applies the cryptos.xlsx data-refresh edit described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.506.02'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '2.184.63'
$ws.Range('E3').Value = '  -1.50%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = "'254.61"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.83%  '
$ws.Range('D6').Value = "'0.607"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.57%  '
$ws.Range('D7').Value = "'73.77"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.46%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -2.75%  '
$ws.Range('D10').Value = "'40.47"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.83%  '
$ws.Range('E11').Value = '  -1.42%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = "'0.101"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'6.78"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.48%  '
$ws.Range('D14').Value = '2.515.90'
$ws.Range('E14').Value = '  -1.41%  '
$ws.Range('D15').Value = "'14.20"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.11%  '
$ws.Range('D16').Value = '2.182.74'
$ws.Range('E16').Value = '  -1.42%  '
$ws.Range('E17').Value = '  -3.50%  '
$ws.Range('D18').Value = '42.424.01'
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('E19').Value = '  -3.01%  '
$ws.Range('D20').Value = "'70.69"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('D22').Value = "'226.96"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.86%  '
$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D23').Value = "'2.13"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.83%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').Value = "'9.35"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -7.04%  '
$ws.Range('E25').Value = '  -0.21%  '
$ws.Range('D26').Value = "'10.50"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.76%  '
$ws.Range('D27').Value = "'3.38"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.23%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = "'2.22"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.76%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').Value = "'2.18"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.82%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = "'36.90"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.36%  '
$ws.Range('D31').Value = "'170.79"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.00%  '
$ws.Range('D32').Value = "'20.04"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.04%  '
$ws.Range('D33').Value = "'0.0809"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.30%  '
$ws.Range('D34').Value = "'5.13"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.63%  '
$ws.Range('E35').Value = '  -1.05%  '
$ws.Range('D36').Value = "'0.107"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.17%  '
$ws.Range('E37').Value = '  -2.79%  '
$ws.Range('E38').Value = '  +5.91%  '
$ws.Range('E39').Value = '  -6.10%  '
$ws.Range('E40').Value = '  -3.45%  '
$ws.Range('E41').Value = '  +0.19%  '
$ws.Range('D42').Value = "'59.35"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.10%  '
$ws.Range('E43').Value = '  -6.58%  '
$ws.Range('D44').Value = "'101.92"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.42%  '
$ws.Range('E45').Value = '  +6.94%  '
$ws.Range('D46').Value = "'0.467"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +9.03%  '
$ws.Range('D47').Value = "'0.0972"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.50%  '
$ws.Range('D48').Value = "'8.22"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.76%  '
$ws.Range('D50').Value = "'1.13"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.82%  '
$ws.Range('E51').Value = '  +0.18%  '
